$wb = $excel.ActiveWorkbook

# Silver Rear_side: B6 5,282 -> 5,263
$ws1 = $wb.Worksheets.Item("Silver Rear_side")
$ws1.Cells.Item(6, 2).NumberFormat = "@"
$ws1.Cells.Item(6, 2).Value = "5,263"

# Silver Busbar front-side: B6 7,907 -> 7,879
$ws2 = $wb.Worksheets.Item("Silver Busbar front-side")
$ws2.Cells.Item(6, 2).NumberFormat = "@"
$ws2.Cells.Item(6, 2).Value = "7,879"

# Silver finger front-side: B6 7,957 -> 7,929
$ws3 = $wb.Worksheets.Item("Silver finger front-side")
$ws3.Cells.Item(6, 2).NumberFormat = "@"
$ws3.Cells.Item(6, 2).Value = "7,929"

# USD_CNY: B6 7.2617 -> 7.2647
$ws4 = $wb.Worksheets.Item("USD_CNY")
$ws4.Cells.Item(6, 2).NumberFormat = "@"
$ws4.Cells.Item(6, 2).Value = "7.2647"
